# Applies the "add Jurisdiction metadata row" edit to the Metadata sheet,
# and updates the Date value, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row at row 11 (pushing Description/Purpose/... down by one),
# so a new "Jurisdiction" property row can be placed right after "Contact".
$ws.Rows.Item(11).Insert()

# New Jurisdiction row - property name + empty value (Value column left blank).
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Copy the style used by all other data rows (s="2") onto the new row's cells.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(11, 2).PasteSpecial(-4122)

# Update the Date value (row 8, column B) to the new timestamp.
$ws.Cells.Item(8, 2).Value = "2024-09-12T14:01:50+00:00"
